# Updates cryptos list values per upstream data refresh (coinranking.com scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.400.48"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.267.06"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.19"
$ws.Range("E5").Value = "  +4.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.32"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.40"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.18"
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.40"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.901"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "2.609.56"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "2.267.16"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "43.552.07"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.84"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.07"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.40"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.86"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.75"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.58"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.20"
$ws.Range("E36").Value = "  +12.03%  "
$ws.Range("E37").Value = "  +8.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.55"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("E40").Value = "  +4.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.86"
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.98"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("E46").Value = "  -7.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.86"
$ws.Range("E47").Value = "  +40.58%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.55"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0995"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.46"
$ws.Range("E51").Value = "  -0.42%  "
